$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3.39
    3  = 2.89
    4  = 2.4
    5  = 2.7
    6  = 3.39
    7  = 3.51
    8  = 3.82
    9  = 4.11
    10 = 4.35
    11 = 4.19
    12 = 4.08
    13 = 4
    14 = 3.97
    15 = 3.98
    16 = 4.07
    17 = 3.79
    18 = 3.78
    19 = 3.6
    20 = 3.62
    21 = 3.7
    22 = 3.63
    23 = 3.66
    24 = 3.59
    25 = 3.41
    26 = 3.94
    27 = 3.55
    28 = 3.43
    29 = 3.33
    30 = 3.08
    31 = 3.29
    32 = 2.84
    33 = 2.51
    34 = 2.63
    35 = 2.72
    36 = 2.87
    37 = 3.08
    38 = 1.77
    39 = 1.48
    40 = 1.43
    41 = 0.96
    42 = 0.75
    43 = 0.61
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
